$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (D) and "Volume(1h)" (E) figures from the latest cryptos refresh.
# Each entry: cell reference, new text value, and whether the text looks like a
# plain number (so it must be forced to Text format to avoid Excel silently
# re-interpreting/truncating it, e.g. "0.07526" -> 7.526E-02).
$updates = @(
    @{ Cell = "D2"; Value = "28.587.11"; IsNumericLooking = $false }
    @{ Cell = "D3"; Value = "1.796.29"; IsNumericLooking = $false }
    @{ Cell = "E3"; Value = "  +0.33%  "; IsNumericLooking = $false }
    @{ Cell = "E4"; Value = "  +0.10%  "; IsNumericLooking = $false }
    @{ Cell = "D5"; Value = "313.65"; IsNumericLooking = $true }
    @{ Cell = "E5"; Value = "  -0.10%  "; IsNumericLooking = $false }
    @{ Cell = "E6"; Value = "  +0.15%  "; IsNumericLooking = $false }
    @{ Cell = "D7"; Value = "0.5278"; IsNumericLooking = $true }
    @{ Cell = "E7"; Value = "  -1.33%  "; IsNumericLooking = $false }
    @{ Cell = "D8"; Value = "0.3771"; IsNumericLooking = $true }
    @{ Cell = "E8"; Value = "  +0.22%  "; IsNumericLooking = $false }
    @{ Cell = "D9"; Value = "42.68"; IsNumericLooking = $true }
    @{ Cell = "E9"; Value = "  -0.60%  "; IsNumericLooking = $false }
    @{ Cell = "D10"; Value = "0.07526"; IsNumericLooking = $true }
    @{ Cell = "E10"; Value = "  +0.29%  "; IsNumericLooking = $false }
    @{ Cell = "E11"; Value = "  +0.37%  "; IsNumericLooking = $false }
    @{ Cell = "D12"; Value = "1.002"; IsNumericLooking = $true }
    @{ Cell = "E12"; Value = "  +0.12%  "; IsNumericLooking = $false }
    @{ Cell = "D13"; Value = "21.09"; IsNumericLooking = $true }
    @{ Cell = "D14"; Value = "7.497"; IsNumericLooking = $true }
    @{ Cell = "E14"; Value = "  +5.90%  "; IsNumericLooking = $false }
    @{ Cell = "D15"; Value = "6.186"; IsNumericLooking = $true }
    @{ Cell = "E15"; Value = "  +0.20%  "; IsNumericLooking = $false }
    @{ Cell = "D16"; Value = "1.796.62"; IsNumericLooking = $false }
    @{ Cell = "E16"; Value = "  +0.46%  "; IsNumericLooking = $false }
    @{ Cell = "D17"; Value = "90.28"; IsNumericLooking = $true }
    @{ Cell = "E17"; Value = "  -0.60%  "; IsNumericLooking = $false }
    @{ Cell = "D18"; Value = "0.00001066"; IsNumericLooking = $true }
    @{ Cell = "E18"; Value = "  -0.22%  "; IsNumericLooking = $false }
    @{ Cell = "D19"; Value = "0.06464"; IsNumericLooking = $true }
    @{ Cell = "E19"; Value = "  -0.56%  "; IsNumericLooking = $false }
    @{ Cell = "E20"; Value = "  +0.14%  "; IsNumericLooking = $false }
    @{ Cell = "D21"; Value = "17.28"; IsNumericLooking = $true }
    @{ Cell = "E21"; Value = "  +2.02%  "; IsNumericLooking = $false }
    @{ Cell = "D22"; Value = "5.926"; IsNumericLooking = $true }
    @{ Cell = "E22"; Value = "  -0.12%  "; IsNumericLooking = $false }
    @{ Cell = "D23"; Value = "28.609.71"; IsNumericLooking = $false }
    @{ Cell = "E23"; Value = "  +3.69%  "; IsNumericLooking = $false }
    @{ Cell = "E24"; Value = "  -0.42%  "; IsNumericLooking = $false }
    @{ Cell = "D25"; Value = "2.094"; IsNumericLooking = $true }
    @{ Cell = "E25"; Value = "  -0.09%  "; IsNumericLooking = $false }
    @{ Cell = "D26"; Value = "160.44"; IsNumericLooking = $true }
    @{ Cell = "E26"; Value = "  +3.30%  "; IsNumericLooking = $false }
    @{ Cell = "E27"; Value = "  +0.18%  "; IsNumericLooking = $false }
    @{ Cell = "D28"; Value = "2.369"; IsNumericLooking = $true }
    @{ Cell = "E28"; Value = "  -0.61%  "; IsNumericLooking = $false }
    @{ Cell = "D29"; Value = "2.003.10"; IsNumericLooking = $false }
    @{ Cell = "E29"; Value = "  +0.45%  "; IsNumericLooking = $false }
    @{ Cell = "D30"; Value = "124.00"; IsNumericLooking = $true }
    @{ Cell = "E30"; Value = "  +1.77%  "; IsNumericLooking = $false }
    @{ Cell = "D31"; Value = "1.116"; IsNumericLooking = $true }
    @{ Cell = "E31"; Value = "  -0.21%  "; IsNumericLooking = $false }
    @{ Cell = "D32"; Value = "0.1025"; IsNumericLooking = $true }
    @{ Cell = "E32"; Value = "  -0.33%  "; IsNumericLooking = $false }
    @{ Cell = "D33"; Value = "5.700"; IsNumericLooking = $true }
    @{ Cell = "E33"; Value = "  +0.45%  "; IsNumericLooking = $false }
    @{ Cell = "D34"; Value = "3.682"; IsNumericLooking = $true }
    @{ Cell = "E34"; Value = "  +1.99%  "; IsNumericLooking = $false }
    @{ Cell = "D35"; Value = "0.2275"; IsNumericLooking = $true }
    @{ Cell = "E35"; Value = "  +9.27%  "; IsNumericLooking = $false }
    @{ Cell = "D36"; Value = "0.06514"; IsNumericLooking = $true }
    @{ Cell = "E36"; Value = "  +8.22%  "; IsNumericLooking = $false }
    @{ Cell = "D37"; Value = "8.908"; IsNumericLooking = $true }
    @{ Cell = "E37"; Value = "  +2.74%  "; IsNumericLooking = $false }
    @{ Cell = "D38"; Value = "0.02313"; IsNumericLooking = $true }
    @{ Cell = "E38"; Value = "  +1.49%  "; IsNumericLooking = $false }
    @{ Cell = "D39"; Value = "5.058"; IsNumericLooking = $true }
    @{ Cell = "E39"; Value = "  +1.43%  "; IsNumericLooking = $false }
    @{ Cell = "E40"; Value = "  +0.38%  "; IsNumericLooking = $false }
    @{ Cell = "D41"; Value = "0.6285"; IsNumericLooking = $true }
    @{ Cell = "E41"; Value = "  +0.65%  "; IsNumericLooking = $false }
    @{ Cell = "D42"; Value = "1.209"; IsNumericLooking = $true }
    @{ Cell = "E42"; Value = "  +5.72%  "; IsNumericLooking = $false }
    @{ Cell = "D43"; Value = "1.001"; IsNumericLooking = $true }
    @{ Cell = "E43"; Value = "  +0.16%  "; IsNumericLooking = $false }
    @{ Cell = "D44"; Value = "1.393"; IsNumericLooking = $true }
    @{ Cell = "E44"; Value = "  -1.41%  "; IsNumericLooking = $false }
    @{ Cell = "D45"; Value = "13.43"; IsNumericLooking = $true }
    @{ Cell = "E45"; Value = "  +0.55%  "; IsNumericLooking = $false }
    @{ Cell = "E46"; Value = "  +0.84%  "; IsNumericLooking = $false }
    @{ Cell = "D47"; Value = "3.665"; IsNumericLooking = $true }
    @{ Cell = "E47"; Value = "  +0.82%  "; IsNumericLooking = $false }
    @{ Cell = "D48"; Value = "126.99"; IsNumericLooking = $true }
    @{ Cell = "E48"; Value = "  +4.56%  "; IsNumericLooking = $false }
    @{ Cell = "E49"; Value = "  +3.09%  "; IsNumericLooking = $false }
    @{ Cell = "D50"; Value = "1.159"; IsNumericLooking = $true }
    @{ Cell = "E50"; Value = "  +2.34%  "; IsNumericLooking = $false }
    @{ Cell = "D51"; Value = "0.06924"; IsNumericLooking = $true }
    @{ Cell = "E51"; Value = "  +2.63%  "; IsNumericLooking = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.IsNumericLooking) {
        # Force text storage so values like "21.09" or "0.07526" keep their
        # original formatting instead of becoming numbers/scientific notation.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
